# Add "Appointment Name" column (C) with data, and add a new row (7)
# with en5 / app5 / Deputy Director, reflecting a "runtime measurement" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C (copy the bold/underline header formatting from B1)
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$ws.Range("C1").Value = "Appointment Name"

# Appointment Name values for existing rows (2-6), matching the order of
# Entity UEN / Appointment UEN already present in column A/B.
$ws.Range("C2").Value = "Director"
$ws.Range("C3").Value = "HR Manager"
$ws.Range("C4").Value = "Tech Lead"
$ws.Range("C5").Value = "Finance Manager"
$ws.Range("C6").Value = "Deputy Director"

# New row of data
$ws.Range("A7").Value = "en5"
$ws.Range("B7").Value = "app5"
$ws.Range("C7").Value = "Deputy Director"

# Match column C width to fit content, like column A/B
$ws.Range("C1").EntireColumn.ColumnWidth = 16.5

# Update selection to reflect new active cell
$ws.Range("D9").Select()
